# "tried and removed correlation table"
# The second table in the document contains a small "site" correlation
# block. One of its rows (site level "1" with an "(empty)" correlation
# value) was removed while keeping the rest of the table intact.

$d = $word.ActiveDocument

$table = $d.Tables.Item(2)

# Find the row whose first cell reads "1" and whose second cell reads
# "(empty)" and delete just that row. Cell.Range.Text carries trailing
# cell-mark / paragraph-mark control characters (and the table uses
# U+2000 padding spaces), so match with StartsWith rather than -eq.
for ($i = $table.Rows.Count; $i -ge 1; $i--) {
    $row = $table.Rows.Item($i)
    $firstCellText = $row.Cells.Item(1).Range.Text
    $secondCellText = $row.Cells.Item(2).Range.Text
    if ($firstCellText.StartsWith("1") -and $secondCellText.StartsWith("(empty)")) {
        $row.Delete()
        break
    }
}
